{"js": "// Update the date line and every \"AxB=\" equation cell in the table to the\n// new values, matching the authoring diff exactly (one-to-one text swap,\n// each original string is unique in the document).\nconst replacements = [\n  [\"2025-10-07 Tuesday\", \"2025-10-08 Wednesday\"],\n  [\"74\u00d733=\", \"56\u00d754=\"],\n  [\"94\u00d758=\", \"52\u00d716=\"],\n  [\"85\u00d765=\", \"15\u00d777=\"],\n  [\"80\u00d780=\", \"69\u00d742=\"],\n  [\"14\u00d735=\", \"96\u00d742=\"],\n  [\"57\u00d737=\", \"12\u00d723=\"],\n  [\"23\u00d743=\", \"31\u00d720=\"],\n  [\"61\u00d749=\", \"61\u00d781=\"],\n  [\"63\u00d714=\", \"25\u00d749=\"],\n  [\"63\u00d720=\", \"16\u00d782=\"],\n  [\"21\u00d728=\", \"33\u00d730=\"],\n  [\"61\u00d771=\", \"11\u00d751=\"],\n  [\"89\u00d717=\", \"47\u00d773=\"],\n  [\"47\u00d753=\", \"57\u00d755=\"],\n  [\"48\u00d772=\", \"28\u00d756=\"],\n  [\"52\u00d769=\", \"60\u00d712=\"],\n  [\"66\u00d722=\", \"59\u00d712=\"],\n  [\"56\u00d780=\", \"42\u00d741=\"],\n  [\"40\u00d723=\", \"83\u00d720=\"],\n  [\"41\u00d779=\", \"70\u00d770=\"],\n  [\"64\u00d769=\", \"84\u00d711=\"],\n  [\"20\u00d713=\", \"46\u00d793=\"],\n  [\"63\u00d741=\", \"55\u00d762=\"],\n  [\"69\u00d796=\", \"98\u00d786=\"],\n  [\"59\u00d799=\", \"96\u00d790=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-07 Tuesday\", \"2025-10-08 Wednesday\"),\n    @(\"74\u00d733=\", \"56\u00d754=\"),\n    @(\"94\u00d758=\", \"52\u00d716=\"),\n    @(\"85\u00d765=\", \"15\u00d777=\"),\n    @(\"80\u00d780=\", \"69\u00d742=\"),\n    @(\"14\u00d735=\", \"96\u00d742=\"),\n    @(\"57\u00d737=\", \"12\u00d723=\"),\n    @(\"23\u00d743=\", \"31\u00d720=\"),\n    @(\"61\u00d749=\", \"61\u00d781=\"),\n    @(\"63\u00d714=\", \"25\u00d749=\"),\n    @(\"63\u00d720=\", \"16\u00d782=\"),\n    @(\"21\u00d728=\", \"33\u00d730=\"),\n    @(\"61\u00d771=\", \"11\u00d751=\"),\n    @(\"89\u00d717=\", \"47\u00d773=\"),\n    @(\"47\u00d753=\", \"57\u00d755=\"),\n    @(\"48\u00d772=\", \"28\u00d756=\"),\n    @(\"52\u00d769=\", \"60\u00d712=\"),\n    @(\"66\u00d722=\", \"59\u00d712=\"),\n    @(\"56\u00d780=\", \"42\u00d741=\"),\n    @(\"40\u00d723=\", \"83\u00d720=\"),\n    @(\"41\u00d779=\", \"70\u00d770=\"),\n    @(\"64\u00d769=\", \"84\u00d711=\"),\n    @(\"20\u00d713=\", \"46\u00d793=\"),\n    @(\"63\u00d741=\", \"55\u00d762=\"),\n    @(\"69\u00d796=\", \"98\u00d786=\"),\n    @(\"59\u00d799=\", \"96\u00d790=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
